$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.302441578384219
$ws.Range("D2").Value = 0.2588745615190362
$ws.Range("E2").Value = 0.1837691009056499
$ws.Range("F2").Value = 5.062241093603916
$ws.Range("G2").Value = 0.002645520021711585
$ws.Range("I2").Value = 1.600655608578592
$ws.Range("J2").Value = 0.1960669283190839
$ws.Range("L2").Value = 2.01447086552696

$ws.Range("B3").Value = 2.197030882498154
$ws.Range("D3").Value = 0.2312334551508854
$ws.Range("E3").Value = 0.1601327085937783
$ws.Range("F3").Value = 5.017848216952501
$ws.Range("G3").Value = 0.002656060729833133
$ws.Range("I3").Value = 1.625447785421905
$ws.Range("J3").Value = 0.1705955069015488
$ws.Range("L3").Value = 1.916170953061169

$ws.Range("B4").Value = 2.133333001872245
$ws.Range("D4").Value = 0.2144484745780915
$ws.Range("E4").Value = 0.1456664053538219
$ws.Range("F4").Value = 4.994673423203551
$ws.Range("G4").Value = 0.002662859040017201
$ws.Range("I4").Value = 1.64182562710743
$ws.Range("J4").Value = 0.1549477588310282
$ws.Range("L4").Value = 1.857005387659171

$ws.Range("B5").Value = 2.107632123492067
$ws.Range("D5").Value = 0.2076531852166283
$ws.Range("E5").Value = 0.139781754203959
$ws.Range("F5").Value = 4.986245066676162
$ws.Range("G5").Value = 0.002665711825277082
$ws.Range("I5").Value = 1.648788605208729
$ws.Range("J5").Value = 0.1485682412089631
$ws.Range("L5").Value = 1.833190334041603

$ws.Range("B6").Value = 2.103379974407119
$ws.Range("D6").Value = 0.2065274692869252
$ws.Range("E6").Value = 0.1388052105848772
$ws.Range("F6").Value = 4.984906562460452
$ws.Range("G6").Value = 0.002666190516756736
$ws.Range("I6").Value = 1.649962201952896
$ws.Range("J6").Value = 0.1475087250207991
$ws.Range("L6").Value = 1.829253586784318

$ws.Range("B7").Value = 2.132985353608547
$ws.Range("D7").Value = 0.2143566527623761
$ws.Range("E7").Value = 0.1455870020268364
$ws.Range("F7").Value = 4.994555657910581
$ws.Range("G7").Value = 0.002662897179751868
$ws.Range("I7").Value = 1.64191836477486
$ws.Range("J7").Value = 0.1548617354374642
$ws.Range("L7").Value = 1.85668301877493

$ws.Range("B8").Value = 2.265882531050806
$ws.Range("D8").Value = 0.2493037768369106
$ws.Range("E8").Value = 0.1756088663528033
$ws.Range("F8").Value = 5.046080352542106
$ws.Range("G8").Value = 0.002649086966389506
$ws.Range("I8").Value = 1.608963165412469
$ws.Range("J8").Value = 0.1872855415949886
$ws.Range("L8").Value = 1.980327335591255

$ws.Range("B9").Value = 2.534696864108696
$ws.Range("D9").Value = 0.3194226651191059
$ws.Range("E9").Value = 0.2349097020296682
$ws.Range("F9").Value = 5.180031380629686
$ws.Range("G9").Value = 0.002624576837461151
$ws.Range("I9").Value = 1.553579357153907
$ws.Range("J9").Value = 0.2508470624670736
$ws.Range("L9").Value = 2.232447692041376

$ws.Range("B10").Value = 2.737327076323538
$ws.Range("D10").Value = 0.3720604277122561
$ws.Range("E10").Value = 0.2788261412495814
$ws.Range("F10").Value = 5.29925441209059
$ws.Range("G10").Value = 0.00260811281998791
$ws.Range("I10").Value = 1.518625090934293
$ws.Range("J10").Value = 0.2975990139934481
$ws.Range("L10").Value = 2.423878387318098

$ws.Range("B11").Value = 2.83065551428615
$ws.Range("D11").Value = 0.3962846334631536
$ws.Range("E11").Value = 0.2989001996021585
$ws.Range("F11").Value = 5.358177446788886
$ws.Range("G11").Value = 0.002600952850820371
$ws.Range("I11").Value = 1.503991670560268
$ws.Range("J11").Value = 0.3188945089902404
$ws.Range("L11").Value = 2.512382857398507

$ws.Range("B12").Value = 2.866164401262836
$ws.Range("D12").Value = 0.4055006315553271
$ws.Range("E12").Value = 0.306517147777484
$ws.Range("F12").Value = 5.38117790932165
$ws.Range("G12").Value = 0.002598288540916501
$ws.Range("I12").Value = 1.498634577234021
$ws.Range("J12").Value = 0.3269637601908926
$ws.Range("L12").Value = 2.546107324408695

$ws.Range("B13").Value = 2.858509448198845
$ws.Range("D13").Value = 0.4035138535046485
$ws.Range("E13").Value = 0.3048759951957578
$ws.Range("F13").Value = 5.376193549263604
$ws.Range("G13").Value = 0.002598860262939284
$ws.Range("I13").Value = 1.49978009548601
$ws.Range("J13").Value = 0.3252256571549879
$ws.Range("L13").Value = 2.538834741848234

$ws.Range("B14").Value = 2.833573484615215
$ws.Range("D14").Value = 0.397041965362547
$ws.Range("E14").Value = 0.2995265343691784
$ws.Range("F14").Value = 5.36005584999657
$ws.Range("G14").Value = 0.002600732716469175
$ws.Range("I14").Value = 1.503547234680624
$ws.Range("J14").Value = 0.3195582622159634
$ws.Range("L14").Value = 2.515153154859036

$ws.Range("B15").Value = 2.818321351385407
$ws.Range("D15").Value = 0.3930834042374443
$ws.Range("E15").Value = 0.2962518778528818
$ws.Range("F15").Value = 5.350260993289652
$ws.Range("G15").Value = 0.002601885759874847
$ws.Range("I15").Value = 1.505878773259852
$ws.Range("J15").Value = 0.3160875183161238
$ws.Range("L15").Value = 2.500674975063248

$ws.Range("B16").Value = 2.731251164848061
$ws.Range("D16").Value = 0.3704831622215465
$ws.Range("E16").Value = 0.2775163180473044
$ws.Range("F16").Value = 5.295499129608288
$ws.Range("G16").Value = 0.002608587340578396
$ws.Range("I16").Value = 1.51960709441925
$ws.Range("J16").Value = 0.2962079567687965
$ws.Range("L16").Value = 2.418123469142813

$ws.Range("B17").Value = 2.678132535377983
$ws.Range("D17").Value = 0.3566920047637154
$ws.Range("E17").Value = 0.2660483614862699
$ws.Range("F17").Value = 5.263114660517431
$ws.Range("G17").Value = 0.002612782690935067
$ws.Range("I17").Value = 1.528355015171158
$ws.Range("J17").Value = 0.2840203871840856
$ws.Range("L17").Value = 2.367848634928691

$ws.Range("B18").Value = 2.647688312074195
$ws.Range("D18").Value = 0.3487857933930911
$ws.Range("E18").Value = 0.259461263301219
$ws.Range("F18").Value = 5.244928389727988
$ws.Range("G18").Value = 0.002615226791293289
$ws.Range("I18").Value = 1.533505763903165
$ws.Range("J18").Value = 0.2770129668527375
$ws.Range("L18").Value = 2.339065366115108

$ws.Range("B19").Value = 2.637398949445185
$ws.Range("D19").Value = 0.346113287459076
$ws.Range("E19").Value = 0.2572324803017665
$ws.Range("F19").Value = 5.238846068137491
$ws.Range("G19").Value = 0.002616059665534388
$ws.Range("I19").Value = 1.535270126963816
$ws.Range("J19").Value = 0.2746407790643843
$ws.Range("L19").Value = 2.3293426137908

$ws.Range("B20").Value = 2.683775891598316
$ws.Range("D20").Value = 0.3581573767237387
$ws.Range("E20").Value = 0.2672682062792688
$ws.Range("F20").Value = 5.266516364288549
$ws.Range("G20").Value = 0.002612332878443127
$ws.Range("I20").Value = 1.52741143539135
$ws.Range("J20").Value = 0.285317503210706
$ws.Range("L20").Value = 2.373186623750712

$ws.Range("B21").Value = 2.840893221381862
$ws.Range("D21").Value = 0.3989417319250776
$ws.Range("E21").Value = 0.3010973722539916
$ws.Range("F21").Value = 5.364777113392506
$ws.Range("G21").Value = 0.002600181458431772
$ws.Range("I21").Value = 1.50243571746843
$ws.Range("J21").Value = 0.3212227657460005
$ws.Range("L21").Value = 2.522103275266488

$ws.Range("B22").Value = 2.944555407318717
$ws.Range("D22").Value = 0.4258475385305189
$ws.Range("E22").Value = 0.3232968180071225
$ws.Range("F22").Value = 5.433010484596224
$ws.Range("G22").Value = 0.002592513664792222
$ws.Range("I22").Value = 1.487187642444447
$ws.Range("J22").Value = 0.3447192591246449
$ws.Range("L22").Value = 2.620654080909901

$ws.Range("B23").Value = 2.889138882953432
$ws.Range("D23").Value = 0.4114635394075776
$ws.Range("E23").Value = 0.3114398019882287
$ws.Range("F23").Value = 5.396221199493368
$ws.Range("G23").Value = 0.002596581178678098
$ws.Range("I23").Value = 1.495226790550085
$ws.Range("J23").Value = 0.3321755846008614
$ws.Range("L23").Value = 2.567941750953594

$ws.Range("B24").Value = 2.681224236408809
$ws.Range("D24").Value = 0.3574948121541013
$ws.Range("E24").Value = 0.2667166959258793
$ws.Range("F24").Value = 5.26497710947109
$ws.Range("G24").Value = 0.002612536138606895
$ws.Range("I24").Value = 1.527837649663773
$ws.Range("J24").Value = 0.2847310790829738
$ws.Range("L24").Value = 2.37077294437438

$ws.Range("B25").Value = 2.461084152938611
$ws.Range("D25").Value = 0.30026802965142
$ws.Range("E25").Value = 0.2188125682297937
$ws.Range("F25").Value = 5.140194690225144
$ws.Range("G25").Value = 0.002630934655616321
$ws.Range("I25").Value = 1.567562332442947
$ws.Range("J25").Value = 0.2336475606307857
$ws.Range("L25").Value = 2.163177272260839
